# Re-export rows 2..12: B (id) becomes "#" + lowercase(C), and D (is_prefered)
# is cleared (no is_pref). Also no lev-distance-based id matching anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = @{
    2  = "Jacora"
    3  = "Griet"
    4  = "Jacoba"
    5  = "Jans"
    6  = "Hans"
    7  = "Geertruy"
    8  = "Andries"
    9  = "Jan"
    10 = "Geertrui"
    11 = "Geererui"
    12 = "Tryn"
}

foreach ($row in 2..12) {
    $name = $newC[$row]
    $ws.Cells.Item($row, 2).Value = "#" + $name.ToLower()
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = ""
}
